$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 (I0) and J1 (IF), matching the existing header style (e.g. H1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# New data values for columns I and J, rows 2-21
$values = @{
    2  = @(1, 4)
    3  = @(1, 6)
    4  = @(1, 5)
    5  = @(1, 5)
    6  = @(4, 7)
    7  = @(9, 9)
    8  = @(7, 9)
    9  = @(8, 8)
    10 = @(2, 6)
    11 = @(1, 4)
    12 = @(1, 6)
    13 = @(6, 7)
    14 = @(8, 8)
    15 = @(1, 6)
    16 = @(1, 6)
    17 = @(1, 5)
    18 = @(1, 5)
    19 = @(1, 7)
    20 = @(1, 2)
    21 = @(7, 8)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
